$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- values that were previously in Row 5
$ws.Range("A2").Value = 111936773
$ws.Range("Q2").Value = 490003.2216792626
$ws.Range("R2").Value = 7087487.484739743
$ws.Range("AC2").Value = "ringhack gamla"

# Row 3 <- values that were previously in Row 4
$ws.Range("A3").Value = 111936769
$ws.Range("Q3").Value = 489837.9886968454
$ws.Range("R3").Value = 7087500.341290037
$ws.Range("AC3").Value = "ringhack gamla"

# Row 4 <- values that were previously in Row 3
$ws.Range("A4").Value = 111936770
$ws.Range("Q4").Value = 489836.5464571039
$ws.Range("R4").Value = 7087463.372540069
$ws.Range("AC4").Value = "ringhack färska"

# Row 5 <- values that were previously in Row 2
$ws.Range("A5").Value = 111936772
$ws.Range("Q5").Value = 489837.4536452024
$ws.Range("R5").Value = 7087471.292509499
$ws.Range("AC5").Value = "ringhack färska"
